$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.319.02'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.651.00'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.09'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.48'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0875'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '1.886.97'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = '1.675.64'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.569'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.43'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = '27.338.48'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.95'
$ws.Range('E18').Value = '  -5.95%  '
$ws.Range('D19').Value = '0.0₃0726'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.40'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  -2.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.42'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.34'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.83'
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.112'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('D33').Value = '1.431.34'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.904'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  +3.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.95'
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.789'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = '1.794.95'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.66'
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.06'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.75'
$ws.Range('E51').Value = '  -0.53%  '
